$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in "Maladie" for the three sick-leave days (rows 12, 13, 14)
$ws.Range("C12").Value = "Maladie"
$ws.Range("C13").Value = "Maladie"
$ws.Range("C14").Value = "Maladie"

# Row 16 (26.02.2018): task + comment about DB export / creation script
$ws.Range("C16").Value = "Correction MLD (envoyé le 08.02 et le 09.02) + création base de données"
$ws.Range("F16").Value = "Exportation de la base de donnée pas terminé + création script création donnée + db"

# Increase row 16 height to fit the wrapped text, and center/wrap its content
$ws.Rows.Item(16).RowHeight = 28.5
$ws.Range("B16").VerticalAlignment = -4108  # xlCenter
$ws.Range("E16").VerticalAlignment = -4108  # xlCenter
$ws.Range("F16:H16").VerticalAlignment = -4108  # xlCenter
$ws.Range("F16:H16").WrapText = $true

# Row 17 (27.02.2018): task + comment about DB creation / data insertion
$ws.Range("C17").Value = "Base de données crée + insertion de données OK + lecture openclasseroom modèle MVC"
$ws.Range("F17").Value = "1'000 données de chaque types ont été crées"

# Update the sheet view: scroll so row 4 is at top, and select F20:H20
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("F20:H20").Select()
